# Auto-generated Excel COM-interop script to apply value updates
# to the Kraken_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 425.46667
$ws.Range("I12").Value = 710.5
$ws.Range("K12").Value = 710.5
$ws.Range("M12").Value = -540.5
$ws.Range("H20").Value = 20000
$ws.Range("I20").Value = 20000
$ws.Range("K20").Value = 20000
$ws.Range("M20").Value = -19770
$ws.Range("H32").Value = 7071.5713
$ws.Range("J32").Value = 9300
$ws.Range("L32").Value = 9300
$ws.Range("N32").Value = -9952
$ws.Range("H35").Value = 20000
$ws.Range("I35").Value = 20000
$ws.Range("K35").Value = 20000
$ws.Range("M35").Value = -19621
$ws.Range("H86").Value = 1000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = -3246
$ws.Range("H89").Value = 1000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = -16232
$ws.Range("H106").Value = 3000
$ws.Range("J106").Value = 3000
$ws.Range("L106").Value = 3000
$ws.Range("N106").Value = -4262
$ws.Range("H130").Value = 96383.8
$ws.Range("J130").Value = 96383.8
$ws.Range("L130").Value = 96383.8
$ws.Range("N130").Value = -106423.8
$ws.Range("H132").Value = 4157.4375
$ws.Range("I132").Value = 4243.6665
$ws.Range("J132").Value = 4046.5715
$ws.Range("K132").Value = 12730.9995
$ws.Range("L132").Value = 12139.7145
$ws.Range("M132").Value = -10200.9995
$ws.Range("N132").Value = -17199.7145
$ws.Range("H135").Value = 3999.6667
$ws.Range("I135").Value = 3999.6667
$ws.Range("K135").Value = 35997.0003
$ws.Range("M135").Value = -33462.0003
$ws.Range("H137").Value = 3253.1428
$ws.Range("I137").Value = 2753.8333
$ws.Range("K137").Value = 8261.499899999999
$ws.Range("M137").Value = -5711.499899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = ""
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = ""
$ws.Range("H121").Value = 99995
$ws.Range("J121").Value = 99995
$ws.Range("L121").Value = 99995
$ws.Range("N121").Value = -103489
$ws.Range("H132").Value = 1995
$ws.Range("I132").Value = 1995
$ws.Range("K132").Value = 5985
$ws.Range("M132").Value = -3455
$ws.Range("H134").Value = 95994
$ws.Range("J134").Value = 95994
$ws.Range("L134").Value = 95994
$ws.Range("N134").Value = -106134

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5255.222
$ws.Range("I94").Value = 4428.143
$ws.Range("J94").Value = 8150
$ws.Range("K94").Value = 4428.143
$ws.Range("L94").Value = 8150
$ws.Range("M94").Value = -3977.143
$ws.Range("N94").Value = -9052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 21018.182
$ws.Range("J50").Value = 19022.223
$ws.Range("L50").Value = 19022.223
$ws.Range("N50").Value = -20272.223
$ws.Range("H51").Value = 17028.572
$ws.Range("J51").Value = 17028.572
$ws.Range("L51").Value = 17028.572
$ws.Range("N51").Value = -18500.572
$ws.Range("H59").Value = 29475
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 29475
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 29475
$ws.Range("M59").Value = ""
$ws.Range("N59").Value = -31765
$ws.Range("H60").Value = 16933.334
$ws.Range("I60").Value = 16500
$ws.Range("J60").Value = 17020
$ws.Range("K60").Value = 16500
$ws.Range("L60").Value = 17020
$ws.Range("M60").Value = -15989
$ws.Range("N60").Value = -18042
$ws.Range("H61").Value = 17028.572
$ws.Range("J61").Value = 17028.572
$ws.Range("L61").Value = 17028.572
$ws.Range("N61").Value = -17724.572
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H86").Value = 4999.5
$ws.Range("I86").Value = 4999.5
$ws.Range("K86").Value = 4999.5
$ws.Range("M86").Value = -3876.5
$ws.Range("H89").Value = 4999.5
$ws.Range("I89").Value = 4999.5
$ws.Range("K89").Value = 24997.5
$ws.Range("M89").Value = -19381.5
$ws.Range("H105").Value = 2015.6666
$ws.Range("I105").Value = 1723.5
$ws.Range("K105").Value = 1723.5
$ws.Range("M105").Value = 23.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 51.0625
$ws.Range("J2").Value = 89.166664
$ws.Range("L2").Value = 534.999984
$ws.Range("N2").Value = -760.999984
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = ""
$ws.Range("H15").Value = 131.66667
$ws.Range("J15").Value = 147.5
$ws.Range("L15").Value = 442.5
$ws.Range("N15").Value = -722.5
$ws.Range("H34").Value = 840.3333
$ws.Range("I34").Value = 197.5
$ws.Range("K34").Value = 592.5
$ws.Range("M34").Value = -508.5
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").Value = ""
$ws.Range("H134").Value = 1200
$ws.Range("J134").Value = 1200
$ws.Range("L134").Value = 3600
$ws.Range("N134").Value = -13740

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").Value = ""
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = ""
$ws.Range("H113").Value = 1749.75
$ws.Range("J113").Value = 1499.5
$ws.Range("L113").Value = 1499.5
$ws.Range("N113").Value = -5839.5
$ws.Range("H119").Value = 99995
$ws.Range("J119").Value = 99995
$ws.Range("L119").Value = 99995
$ws.Range("N119").Value = -109671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2997.3333
$ws.Range("I7").Value = 1596.8
$ws.Range("K7").Value = 1596.8
$ws.Range("M7").Value = -1484.8
$ws.Range("H46").Value = 3159.8235
$ws.Range("I46").Value = 2233.6667
$ws.Range("J46").Value = 3358.2856
$ws.Range("K46").Value = 2233.6667
$ws.Range("L46").Value = 3358.2856
$ws.Range("M46").Value = -2045.6667
$ws.Range("N46").Value = -3734.2856
$ws.Range("H55").Value = 1718
$ws.Range("J55").Value = 2465
$ws.Range("L55").Value = 2465
$ws.Range("N55").Value = -2811
$ws.Range("H100").Value = 7242.1665
$ws.Range("I100").Value = 7242.1665
$ws.Range("K100").Value = 7242.1665
$ws.Range("M100").Value = -6701.1665
$ws.Range("H126").Value = 2997.3333
$ws.Range("I126").Value = 1596.8
$ws.Range("K126").Value = 4790.4
$ws.Range("M126").Value = -2320.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 17170.5
$ws.Range("J95").Value = 17170.5
$ws.Range("L95").Value = 17170.5
$ws.Range("N95").Value = -22662.5
$ws.Range("H118").Value = 118742.5
$ws.Range("I118").Value = 99990
$ws.Range("K118").Value = 99990
$ws.Range("M118").Value = -98333
$ws.Range("H132").Value = 5113
$ws.Range("J132").Value = 7880
$ws.Range("L132").Value = 23640
$ws.Range("N132").Value = -28700
